# Header updates for summer uploads
# Rename several header cells in row 1 of the (single) worksheet to match
# the new Ministry reporting header labels, then leave the header row
# (A1:K1) selected - mirroring the authored edit's sheetView selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

$ws.Range("A1:K1").Select()
